# Apply the cryptos-list refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "'30.778.65"
$ws.Range("E2").Value = "  -1.06%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "'1.941.69"
$ws.Range("E3").Value = "  -0.90%  "

# Row 4: TetherUSD
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.05%  "

# Row 5: BNB
$ws.Range("D5").Value = "'242.02"
$ws.Range("E5").Value = "  -1.84%  "

# Row 6: USDC
$ws.Range("E6").Value = "  -0.11%  "

# Row 7: XRP
$ws.Range("D7").Value = "'0.4894"
$ws.Range("E7").Value = "  +0.05%  "

# Row 8: Cardano
$ws.Range("D8").Value = "'0.2956"
$ws.Range("E8").Value = "  -0.40%  "

# Row 9: Dogecoin
$ws.Range("D9").Value = "'0.06894"
$ws.Range("E9").Value = "  +1.11%  "

# Row 10: Solana
$ws.Range("D10").Value = "'19.43"
$ws.Range("E10").Value = "  +2.00%  "

# Row 11: Litecoin
$ws.Range("D11").Value = "'106.59"
$ws.Range("E11").Value = "  +0.00%  "

# Row 12: WrappedEther
$ws.Range("D12").Value = "'1.939.69"
$ws.Range("E12").Value = "  -0.73%  "

# Row 13: TRON
$ws.Range("D13").Value = "'0.07717"
$ws.Range("E13").Value = "  -0.45%  "

# Row 14: Polkadot
$ws.Range("D14").Value = "'5.345"
$ws.Range("E14").Value = "  -1.22%  "

# Row 15: Polygon
$ws.Range("D15").Value = "'0.6988"
$ws.Range("E15").Value = "  -1.52%  "

# Row 16: BitcoinCash
$ws.Range("D16").Value = "'277.53"
$ws.Range("E16").Value = "  -2.07%  "

# Row 17: WrappedBTC
$ws.Range("D17").Value = "'30.801.24"
$ws.Range("E17").Value = "  -0.62%  "

# Row 18: WrappedliquidstakedEther2.0
$ws.Range("B18").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C18").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D18").Value = "'2.233.78"
$ws.Range("E18").Value = "  +1.87%  "

# Row 19: ShibaInu
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000007718"
$ws.Range("E19").Value = "  -0.45%  "

# Row 20: Avalanche
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "'13.10"
$ws.Range("E20").Value = "  -0.83%  "

# Row 21: Dai
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'0.9998"
$ws.Range("E21").Value = "  -0.20%  "

# Row 22: Uniswap
$ws.Range("D22").Value = "'5.468"
$ws.Range("E22").Value = "  -1.44%  "

# Row 23: BinanceUSD
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.07%  "

# Row 24: Chainlink
$ws.Range("D24").Value = "'6.519"
$ws.Range("E24").Value = "  -1.08%  "

# Row 25: Cosmos
$ws.Range("D25").Value = "'9.720"
$ws.Range("E25").Value = "  -2.39%  "

# Row 26: Monero
$ws.Range("D26").Value = "'167.66"
$ws.Range("E26").Value = "  -0.65%  "

# Row 27: EthereumClassic
$ws.Range("D27").Value = "'19.64"
$ws.Range("E27").Value = "  -1.78%  "

# Row 28: LidoDAOToken
$ws.Range("D28").Value = "'2.152"
$ws.Range("E28").Value = "  -1.81%  "

# Row 29: Stellar
$ws.Range("E29").Value = "  -1.11%  "

# Row 30: Toncoin
$ws.Range("D30").Value = "'1.387"
$ws.Range("E30").Value = "  -3.78%  "

# Row 31: PancakeSwap
$ws.Range("D31").Value = "'1.550"
$ws.Range("E31").Value = "  -2.67%  "

# Row 32: Filecoin
$ws.Range("D32").Value = "'4.548"
$ws.Range("E32").Value = "  -5.14%  "

# Row 33: InternetComputer(DFINITY)
$ws.Range("D33").Value = "'4.357"
$ws.Range("E33").Value = "  -3.50%  "

# Row 34: Hedera
$ws.Range("D34").Value = "'0.04849"
$ws.Range("E34").Value = "  -3.05%  "

# Row 35: ImmutableX
$ws.Range("D35").Value = "'0.7524"
$ws.Range("E35").Value = "  -1.85%  "

# Row 36: ARBITRUM
$ws.Range("D36").Value = "'1.159"
$ws.Range("E36").Value = "  -0.57%  "

# Row 37: Frax
$ws.Range("E37").Value = "  -0.11%  "

# Row 38: HuobiToken
$ws.Range("D38").Value = "'2.726"
$ws.Range("E38").Value = "  -0.34%  "

# Row 39: VeChain
$ws.Range("D39").Value = "'0.01993"
$ws.Range("E39").Value = "  -2.72%  "

# Row 40: MXToken
$ws.Range("E40").Value = "  -1.73%  "

# Row 41: Aave
$ws.Range("D41").Value = "'78.07"
$ws.Range("E41").Value = "  +6.42%  "

# Row 42: FraxShare
$ws.Range("D42").Value = "'6.483"
$ws.Range("E42").Value = "  +0.89%  "

# Row 43: RenderToken
$ws.Range("E43").Value = "  -1.89%  "

# Row 44: TrustWalletToken
$ws.Range("D44").Value = "'0.9056"
$ws.Range("E44").Value = "  +2.64%  "

# Row 45: Quant
$ws.Range("D45").Value = "'108.03"
$ws.Range("E45").Value = "  -1.13%  "

# Row 46: TheSandbox
$ws.Range("D46").Value = "'0.4397"
$ws.Range("E46").Value = "  -1.20%  "

# Row 47: PaxDollar
$ws.Range("D47").Value = "'0.9984"
$ws.Range("E47").Value = "  -0.25%  "

# Row 48: Aptos
$ws.Range("D48").Value = "'7.739"
$ws.Range("E48").Value = "  +3.67%  "

# Row 49: Maker
$ws.Range("D49").Value = "'987.58"
$ws.Range("E49").Value = "  -2.91%  "

# Row 51: EnergySwap
$ws.Range("D51").Value = "'9.286"
$ws.Range("E51").Value = "  -0.84%  "
